$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E ("purpose") rows 2-25 contain "fullRNASEQ" and should read "fullRNASeq"
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 5).Value = "fullRNASeq"
}
